$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 74, shifting rows 74:144 down to 75:145
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new data
$ws.Range("A74").Value = 8
$ws.Range("B74").Value = "Terminal La Palmera de La Serena"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = 44705
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
$ws.Range("E74").Value = 4
$ws.Range("F74").Value = 100112044
$ws.Range("G74").Value = "Perejil"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 2400
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = 1750
$ws.Range("N74").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O74").Value = "Provincia del Elquí"
$ws.Range("P74").Value = 1167
$ws.Range("Q74").Value = 1.5
$ws.Range("R74").Value = "Hortaliza"
